# Add a new "2022-Q3" sheet (with its fund-holding detail) and a
# corresponding summary row on the "总计" sheet, pushing the existing
# quarters down by one row/position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet with the same look & feel
#    as the existing per-quarter sheets (copy formatting from the
#    sheet that is currently "2022-Q2" - soon to become the neighbour
#    of our new sheet).
# ---------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

# Copy the header-row formatting (bold, border, centered) ...
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ... and the formatting used for the small numeric index in column A.
$refSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header values
$newSheet.Cells.Item(1, 2).Value2 = "基金代码"
$newSheet.Cells.Item(1, 3).Value2 = "基金名称"
$newSheet.Cells.Item(1, 4).Value2 = "基金规模"
$newSheet.Cells.Item(1, 5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value2 = "仓位占比"
$newSheet.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value2 = "仓位排名"

# Data row (one fund holding reported for 2022-Q3)
$newSheet.Cells.Item(2, 1).Value2 = 0
$newSheet.Cells.Item(2, 2).Value2 = "'005281"
$newSheet.Cells.Item(2, 3).Value2 = "中科沃土转型升级灵活配置混合"
$newSheet.Cells.Item(2, 4).Value2 = "'0.10"
$newSheet.Cells.Item(2, 5).Value2 = "'57.70"
$newSheet.Cells.Item(2, 6).Value2 = "'2.79"
$newSheet.Cells.Item(2, 7).Value2 = "'0.0028"
$newSheet.Cells.Item(2, 8).Value2 = 10

# Move the brand-new sheet so it sits right before "2022-Q2", i.e.
# immediately after the "总计" summary sheet.
$newSheet.Move($refSheet)

# ---------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the six existing quarterly
#    rows down by one and insert the "2022-Q3" totals at the top.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 7; $r -ge 2; $r--) {
    $newR = $r + 1
    $summary.Cells.Item($newR, 2).Value2 = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($newR, 3).Value2 = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($newR, 4).Value2 = $summary.Cells.Item($r, 4).Value2
}

$summary.Cells.Item(2, 2).Value2 = "2022-Q3"
$summary.Cells.Item(2, 3).Value2 = 1
$summary.Cells.Item(2, 4).Value2 = 0

# Re-number the column-A running index (0..6) across all eight rows.
for ($r = 2; $r -le 8; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}

# Row 8 is brand new - give its index cell the same formatting as the
# rest of column A by copying it from the row above.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)
